$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (row 1 is the header, rows 2-4 are the records).
# Column A = firstName, B = middleName, C = lastName,
# D = photograph (unchanged), E = username.

$ws.Range("A2").Value = "Dan"
$ws.Range("B2").Value = "Ivan"
$ws.Range("C2").Value = "Don"
$ws.Range("E2").Value = "donc1234"

$ws.Range("A3").Value = "Serg"
$ws.Range("B3").Value = "Ivan"
$ws.Range("C3").Value = "Don"
$ws.Range("E3").Value = "donc4564"

$ws.Range("A4").Value = "Sam"
$ws.Range("B4").Value = "Ivan"
$ws.Range("C4").Value = "Don"
$ws.Range("E4").Value = "donc7894"

# Update the active cell selection to A2.
$ws.Range("A2").Select()
